$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting the existing D:K data (and the
# newly vacated column) one column to the right (now D:L).
$ws.Columns("D:D").Insert()

# The freshly inserted column has no number format yet; copy the format
# from the (now-shifted) column E so each row's new D cell matches the
# style used by the rest of that row (date format for header rows, the
# numeric format elsewhere).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reported quarter (period
# ending 2018-10-31, serial 43404) for each statement, and its data.

# Income Statement
$ws.Range("D7").Value = 43404
$ws.Range("D8").Value = 57300
$ws.Range("D9").Value = 22700
$ws.Range("D10").Value = 34600
$ws.Range("D12").Value = 31700
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 66700
$ws.Range("D18").Value = -9400
$ws.Range("D20").Value = 1000
$ws.Range("D21").Value = -6600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -8400
$ws.Range("D24").Value = 600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -9000
$ws.Range("D27").Value = -9000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1000
$ws.Range("D33").Value = -9000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -9000

# Balance Sheet
$ws.Range("D38").Value = 43404
$ws.Range("D41").Value = 212300
$ws.Range("D42").Value = 136300
$ws.Range("D43").Value = 32400
$ws.Range("D44").Value = 23300
$ws.Range("D45").Value = 3300
$ws.Range("D46").Value = 407700
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 6700
$ws.Range("D49").Value = 37800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 6400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 458600
$ws.Range("D57").Value = 11100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 26800
$ws.Range("D60").Value = 37900
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 9600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 47500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 235400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 411100
$ws.Range("D77").Value = 0

# Cash Flow Statement
$ws.Range("D80").Value = 43404
$ws.Range("D81").Value = -9000
$ws.Range("D83").Value = 1800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 3500
$ws.Range("D91").Value = -400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -55900
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -30600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -82900

Write-Host "Financials updated"
